$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("統計")

$newRow = 25

$ws.Cells.Item($newRow, 1).Value = "2025-08-30T06:22:38.182733"
$ws.Cells.Item($newRow, 2).Value = 6
$ws.Cells.Item($newRow, 3).Value = "全案件リスト"
$ws.Cells.Item($newRow, 4).Value = 66.7
$ws.Cells.Item($newRow, 5).Value = 2
$ws.Cells.Item($newRow, 6).Value = 3
$ws.Cells.Item($newRow, 7).Value = 6
